$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 7325.28
$ws.Range("I19").Value = 3401.2727
$ws.Range("J19").Value = 10408.429
$ws.Range("K19").Value = 3401.2727
$ws.Range("L19").Value = 10408.429
$ws.Range("M19").Value = -3226.2727
$ws.Range("N19").Value = -10758.429

$ws.Range("H33").Value = 485.67648
$ws.Range("I33").Value = 536.4286
$ws.Range("K33").Value = 536.4286
$ws.Range("M33").Value = -307.4286

$ws.Range("H131").Value = 6353.364
$ws.Range("I131").Value = 4269.5713
$ws.Range("J131").Value = 10000
$ws.Range("K131").Value = 12808.7139
$ws.Range("L131").Value = 30000
$ws.Range("M131").Value = -7768.713899999999
$ws.Range("N131").Value = -40080

$ws.Range("H134").Value = 152499.5
$ws.Range("J134").Value = 152499.5
$ws.Range("L134").Value = 152499.5
$ws.Range("N134").Value = -162639.5

$ws.Range("H137").Value = 4874.1797
$ws.Range("I137").Value = 6323.4165
$ws.Range("K137").Value = 18970.2495
$ws.Range("M137").Value = -16420.2495

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5751.5
$ws.Range("I32").Value = 4666.288
$ws.Range("K32").Value = 4666.288
$ws.Range("M32").Value = -4379.288

$ws.Range("H61").Value = 3169.4814
$ws.Range("J61").Value = 3999
$ws.Range("L61").Value = 3999
$ws.Range("N61").Value = -4423

$ws.Range("H74").Value = 3099.0908
$ws.Range("I74").Value = 2609
$ws.Range("K74").Value = 2609
$ws.Range("M74").Value = -1735

$ws.Range("H77").Value = 3099.0908
$ws.Range("I77").Value = 2609
$ws.Range("K77").Value = 13045
$ws.Range("M77").Value = -8677

$ws.Range("H136").Value = 3169.4814
$ws.Range("J136").Value = 3999
$ws.Range("L136").Value = 11997
$ws.Range("N136").Value = -17097

$ws.Range("H141").Value = 145326.86
$ws.Range("J141").Value = 145326.86
$ws.Range("L141").Value = 145326.86
$ws.Range("N141").Value = -155686.86

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 19608364
$ws.Range("I64").Value = 33333784
$ws.Range("J64").Value = 620.8570999999999
$ws.Range("K64").Value = 33333784
$ws.Range("L64").Value = 620.8570999999999
$ws.Range("M64").Value = -33333559
$ws.Range("N64").Value = -1070.8571

$ws.Range("H67").Value = 19608364
$ws.Range("I67").Value = 33333784
$ws.Range("J67").Value = 620.8570999999999
$ws.Range("K67").Value = 33333784
$ws.Range("L67").Value = 620.8570999999999
$ws.Range("M67").Value = -33333004
$ws.Range("N67").Value = -2180.8571

$ws.Range("H134").Value = 6150.3887
$ws.Range("J134").Value = 6004.75
$ws.Range("L134").Value = 18014.25
$ws.Range("N134").Value = -23084.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5999
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H31").Value = 1892.7188
$ws.Range("I31").Value = 1352.037
$ws.Range("K31").Value = 1352.037
$ws.Range("M31").Value = -1057.037

$ws.Range("H34").Value = 1892.7188
$ws.Range("I34").Value = 1352.037
$ws.Range("K34").Value = 1352.037
$ws.Range("M34").Value = -1150.037

$ws.Range("H132").Value = 6445.8706
$ws.Range("I132").Value = 3363.946
$ws.Range("K132").Value = 10091.838
$ws.Range("M132").Value = -7561.838

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5301.625
$ws.Range("I3").Value = 4702
$ws.Range("J3").Value = 9499
$ws.Range("K3").Value = 14106
$ws.Range("L3").Value = 28497
$ws.Range("M3").Value = -13994
$ws.Range("N3").Value = -28721

$ws.Range("H5").Value = 1279.4445
$ws.Range("I5").Value = 502.83334
$ws.Range("K5").Value = 1508.50002
$ws.Range("M5").Value = -1396.50002

$ws.Range("H132").Value = 13999.2
$ws.Range("J132").Value = 13999.2
$ws.Range("L132").Value = 125992.8
$ws.Range("N132").Value = -131052.8

$ws.Range("H135").Value = 1279.4445
$ws.Range("I135").Value = 502.83334
$ws.Range("K135").Value = 4525.50006
$ws.Range("M135").Value = -1990.50006

$ws.Range("H141").Value = 12336.182
$ws.Range("I141").Value = 13744.333
$ws.Range("K141").Value = 41232.999
$ws.Range("M141").Value = -36052.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 68.59999999999999
$ws.Range("J2").Value = 134
$ws.Range("L2").Value = 134
$ws.Range("N2").Value = -360

$ws.Range("H70").Value = 4964.736
$ws.Range("I70").Value = 4907.174
$ws.Range("K70").Value = 4907.174
$ws.Range("M70").Value = -4637.174

$ws.Range("H73").Value = 4964.736
$ws.Range("I73").Value = 4907.174
$ws.Range("K73").Value = 4907.174
$ws.Range("M73").Value = -3971.174

$ws.Range("H80").Value = 34288276
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 34288276
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H132").Value = 5151.0566
$ws.Range("I132").Value = 2901.1738
$ws.Range("K132").Value = 8703.5214
$ws.Range("M132").Value = -6173.5214

$ws.Range("H140").Value = 71943.5
$ws.Range("J140").Value = 82332.2
$ws.Range("L140").Value = 82332.2
$ws.Range("N140").Value = -92692.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4172447.2
$ws.Range("I16").Value = 5686537
$ws.Range("J16").Value = 8700.125
$ws.Range("K16").Value = 5686537
$ws.Range("L16").Value = 8700.125
$ws.Range("M16").Value = -5686367
$ws.Range("N16").Value = -9040.125

$ws.Range("H22").Value = 973
$ws.Range("I22").Value = 965.5454999999999
$ws.Range("J22").Value = 993.5
$ws.Range("K22").Value = 965.5454999999999
$ws.Range("L22").Value = 993.5
$ws.Range("M22").Value = -670.5454999999999
$ws.Range("N22").Value = -1583.5

$ws.Range("H27").Value = 973
$ws.Range("I27").Value = 965.5454999999999
$ws.Range("J27").Value = 993.5
$ws.Range("K27").Value = 965.5454999999999
$ws.Range("L27").Value = 993.5
$ws.Range("M27").Value = -858.5454999999999
$ws.Range("N27").Value = -1207.5

$ws.Range("H42").Value = 21874
$ws.Range("I42").Value = 18748.5
$ws.Range("K42").Value = 18748.5
$ws.Range("M42").Value = -18185.5

$ws.Range("H46").Value = 3180.625
$ws.Range("I46").Value = 2337.125
$ws.Range("J46").Value = 3602.375
$ws.Range("K46").Value = 2337.125
$ws.Range("L46").Value = 3602.375
$ws.Range("M46").Value = -2149.125
$ws.Range("N46").Value = -3978.375

$ws.Range("H49").Value = 21874
$ws.Range("I49").Value = 18748.5
$ws.Range("K49").Value = 18748.5
$ws.Range("M49").Value = -18601.5

$ws.Range("H132").Value = 164179.5
$ws.Range("I132").Value = 232429.72
$ws.Range("J132").Value = 4929
$ws.Range("K132").Value = 697289.16
$ws.Range("L132").Value = 14787
$ws.Range("M132").Value = -694759.16
$ws.Range("N132").Value = -19847

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 36599.5
$ws.Range("J122").Value = 37159.2
$ws.Range("L122").Value = 111477.6
$ws.Range("N122").Value = -116377.6

$ws.Range("H132").Value = 5053.289
$ws.Range("I132").Value = 3605.3713
$ws.Range("K132").Value = 10816.1139
$ws.Range("M132").Value = -8286.1139

$ws.Range("H141").Value = 192000
$ws.Range("J141").Value = 192000
$ws.Range("L141").Value = 192000
$ws.Range("N141").Value = -202360
